$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the three "extra" test rows (tom / ron / " or 1=1 -- -") that were
# added for testing — rows 3-5 on the invalid_logins sheet. This shifts the
# remaining rows (the admin' / admin") SQL-injection rows) up from 6-7 to 3-4,
# matching the trimmed-down test data set, and Excel prunes the now-unused
# shared strings automatically.
$ws.Rows("3:5").Delete()

# Restore the visible selection over the remaining injected rows.
$ws.Range("A3:B4").Select()
